# Change Background Music && Add Login feature
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the merged "Loc Le" record (was row 3's data, with coin added)
$ws.Range("A2").Value = "Loc Le"
$ws.Range("B2").Value = "ltloc05lumia520@gmail.com"
$ws.Range("C2").Value = "Loc Le.png"
$ws.Range("D2").Value = "khongbiethehe"
$ws.Range("E2").Value = 0

# E1 header "coin" becomes numeric 0
$ws.Range("E1").Value = 0

# Drop the old row 3 entirely (data now lives on row 2)
$ws.Range("A3:F3").Delete()

# Update the active selection to match the authored state
$ws.Range("G5").Select()
